# Build site at 2022-09-26 16:07:08 UTC
# LOT2038.xlsx update:
#  - Row 13 (docente name in B13/C13, no label in A13) is deleted entirely,
#    shifting every row below it up by one.
#  - Several long placeholder paragraphs are swapped out for short values
#    (re-using text that appears elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete row 13 (only contained B13/C13 = docente name, no A13 label).
#    Everything below shifts up by one row.
$ws.Rows(13).Delete()

# 2) Objetivos: (row 10) long description -> docente name
$ws.Range("B10").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("C10").Value = "1097178 - João Batista de Almeida e Silva"

# 3) Programa resumido: (row 13 after shift) long description -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 4) Programa: (row 15 after shift) long description -> activation date.
#    Use Copy so the value is carried over as the same shared text string
#    (matching the "Ativação:" row) instead of being auto-parsed as a date.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# 5) Método: (row 18 after shift) -> docente name
$ws.Range("B18").Value = "1097178 - João Batista de Almeida e Silva"
$ws.Range("C18").Value = "1097178 - João Batista de Almeida e Silva"

# 6) Critério: (row 19 after shift) -> old Método text
$ws.Range("B19").Value = "Avaliação teórica, com provas escritas"
$ws.Range("C19").Value = "Avaliação teórica, com provas escritas"

# 7) Norma de recuperação: (row 20 after shift) -> old Critério text
$ws.Range("B20").Value = "Provas e média aritimética das provas"
$ws.Range("C20").Value = "Provas e média aritimética das provas"

# 8) Bibliografia: (row 21 after shift) -> old Norma de recuperação text
$ws.Range("B21").Value = "A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculadapela equação: MF = (NF + PR)/2."
$ws.Range("C21").Value = "A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculadapela equação: MF = (NF + PR)/2."
